{"js": "// 1) \"Gustavo Emerick da Silva\" \u2014 collapse the spell-check-split runs\n//    (Gustavo / Emerick / da Silva) into a single run of plain text.\n//    This occurs twice: once in the participants table, once in the\n//    \"Equipe de Projeto\" paragraph.\n{\n  const results = context.document.body.search(\"Gustavo Emerick da Silva\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(\"Gustavo Emerick da Silva\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// 2) \"Disponibilidade Oportuna: Em situa\u00e7\u00f5es de emerg\u00eancia, ...\" \u2014 merge\n//    the runs that were split apart by a grammar-check marker back into\n//    a single run, keeping the preceding bold label run untouched.\n{\n  const target =\n    \" Em situa\u00e7\u00f5es de emerg\u00eancia, como acidentes graves ou desastres naturais, a disponibilidade oportuna de sangue \u00e9 cr\u00edtica. O projeto visa melhorar a capacidade de resposta em situa\u00e7\u00f5es de alta demanda. \";\n  const results = context.document.body.search(target, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(target, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// 3) Remove the \"Diagrama de Neg\u00f3cios\" heading, its surrounding blank\n//    paragraphs and the business-diagram picture \u2014 everything after the\n//    \"Essas limita\u00e7\u00f5es...\" paragraph up to (and including) the picture's\n//    paragraph, right before the section break.\n{\n  const paragraphs = context.document.body.paragraphs;\n  paragraphs.load(\"items/text\");\n  await context.sync();\n\n  const items = paragraphs.items;\n  let startIndex = -1;\n  let endIndex = -1;\n\n  for (let i = 0; i < items.length; i++) {\n    if (items[i].text.indexOf(\"Essas limita\u00e7\u00f5es destacam os desafios\") !== -1) {\n      startIndex = i + 1; // first paragraph to delete (right after this one)\n      break;\n    }\n  }\n\n  if (startIndex !== -1) {\n    for (let i = startIndex; i < items.length; i++) {\n      if (items[i].text.indexOf(\"Diagrama de Neg\u00f3cios\") !== -1) {\n        // the section continues a couple more (blank) paragraphs plus the\n        // picture paragraph \u2014 delete through the end of the body.\n        endIndex = items.length - 1;\n        break;\n      }\n    }\n  }\n\n  if (startIndex !== -1 && endIndex !== -1) {\n    for (let i = endIndex; i >= startIndex; i--) {\n      items[i].delete();\n    }\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"Gustavo Emerick da Silva\" \u2014 collapse the spell-check-split runs\n#    (Gustavo / Emerick / da Silva) into a single run of plain text.\n#    Find/Replace naturally merges the runs and drops the now-stale\n#    <w:proofErr/> spell-check markers. Occurs twice: once in the\n#    participants table, once in the \"Equipe de Projeto\" paragraph.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Gustavo Emerick da Silva\"\n$find.Replacement.Text = \"Gustavo Emerick da Silva\"\n$find.Execute(\n  [ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false,\n  [ref]$find.Replacement.Text, 2\n) | Out-Null\n\n# 2) \"Disponibilidade Oportuna: Em situa\u00e7\u00f5es de emerg\u00eancia, ...\" \u2014 merge\n#    the runs that a grammar-check marker split apart back into a single\n#    run, leaving the preceding bold label run untouched.\n$target = \" Em situa\u00e7\u00f5es de emerg\u00eancia, como acidentes graves ou desastres naturais, a disponibilidade oportuna de sangue \u00e9 cr\u00edtica. O projeto visa melhorar a capacidade de resposta em situa\u00e7\u00f5es de alta demanda. \"\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = $target\n$find2.Replacement.Text = $target\n$find2.Execute(\n  [ref]$find2.Text, $false, $false, $false, $false, $false, $true, 1, $false,\n  [ref]$find2.Replacement.Text, 2\n) | Out-Null\n\n# 3) Remove the \"Diagrama de Neg\u00f3cios\" heading, its surrounding blank\n#    paragraphs and the business-diagram picture \u2014 everything after the\n#    \"Essas limita\u00e7\u00f5es...\" paragraph through to the end of the body\n#    (right before the section break).\n$count = $d.Paragraphs.Count\n$startIndex = -1\n$endIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n  $text = $d.Paragraphs.Item($i).Range.Text\n  if ($text -like \"*Essas limita*es destacam os desafios*\") {\n    $startIndex = $i + 1\n  }\n  if ($startIndex -ne -1 -and $i -ge $startIndex -and $text -like \"*Diagrama de Neg*cios*\") {\n    $endIndex = $count\n    break\n  }\n}\n\nif ($startIndex -ne -1 -and $endIndex -ne -1) {\n  $startPara = $d.Paragraphs.Item($startIndex)\n  $endPara = $d.Paragraphs.Item($endIndex)\n  $deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)\n  $deleteRange.Delete()\n}\n"}
